$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "NoFluff"

# Clear all existing content/formatting before rebuilding
$ws.Cells.Clear()

# Header row
$ws.Cells.Item(1,1).Value() = 'LINK'
$ws.Cells.Item(1,2).Value() = 'OPIS'
$ws.Cells.Item(1,3).Value() = 'FIRMA'
$ws.Cells.Item(1,4).Value() = 'ZAROBKI'
$ws.Cells.Item(1,5).Value() = 'LOKALIZACJA'

# Header formatting: bold + centered (format A1 directly, then propagate via
# copy/paste-special so the whole header row shares a single cell style)
$a1 = $ws.Range("A1")
$a1.Font.Bold() = $true
$a1.HorizontalAlignment() = -4108
$a1.Copy()
$ws.Range("B1:E1").PasteSpecial(-4122)
$excel.CutCopyMode() = $false

# Job listing rows
$ws.Cells.Item(2,1).Formula() = '=HYPERLINK("https://nofluffjobs.com/pl/job/senior-mobile-qa-automation-engineer-beekeeper-remote", "https://nofluffjobs.com/pl/job/senior-mobile-qa-automation-engineer-beekeeper-remote")'
$ws.Cells.Item(2,2).Value() = ' Senior Mobile QA Automation Engineer '
$ws.Cells.Item(2,3).Value() = ' Beekeeper AG '
$ws.Cells.Item(2,4).Value() = ' 17500  – 26666  PLN'
$ws.Cells.Item(2,5).Value() = ' Zdalnie   +1 '

$ws.Cells.Item(3,1).Formula() = '=HYPERLINK("https://nofluffjobs.com/pl/job/qa-engineer-grape-up-remote-2", "https://nofluffjobs.com/pl/job/qa-engineer-grape-up-remote-2")'
$ws.Cells.Item(3,2).Value() = ' QA Engineer '
$ws.Cells.Item(3,3).Value() = ' Grape Up '
$ws.Cells.Item(3,4).Value() = ' 12000  – 18000  PLN'
$ws.Cells.Item(3,5).Value() = ' Zdalnie   +3 '

$ws.Cells.Item(4,1).Formula() = '=HYPERLINK("https://nofluffjobs.com/pl/job/senior-test-engineer-rpo-partners-gdansk", "https://nofluffjobs.com/pl/job/senior-test-engineer-rpo-partners-gdansk")'
$ws.Cells.Item(4,2).Value() = ' Senior Test Engineer '
$ws.Cells.Item(4,3).Value() = ' RPO Partners Sp. z o.o. '
$ws.Cells.Item(4,4).Value() = ' 12000  – 23000  PLN'
$ws.Cells.Item(4,5).Value() = ' Gdańsk  '

$ws.Cells.Item(5,1).Formula() = '=HYPERLINK("https://nofluffjobs.com/pl/job/test-automation-engineer-sii-polska-remote-7", "https://nofluffjobs.com/pl/job/test-automation-engineer-sii-polska-remote-7")'
$ws.Cells.Item(5,2).Value() = ' Test Automation Engineer '
$ws.Cells.Item(5,3).Value() = ' Sii Polska '
$ws.Cells.Item(5,4).Value() = ' 21840  – 25410  PLN'
$ws.Cells.Item(5,5).Value() = ' Zdalnie   +5 '

$ws.Cells.Item(6,1).Formula() = '=HYPERLINK("https://nofluffjobs.com/pl/job/tester-automatyzujacy-remote-link-group", "https://nofluffjobs.com/pl/job/tester-automatyzujacy-remote-link-group")'
$ws.Cells.Item(6,2).Value() = ' Tester Automatyzujący - Remote '
$ws.Cells.Item(6,3).Value() = ' Link Group '
$ws.Cells.Item(6,4).Value() = ' 16800  – 21840  PLN'
$ws.Cells.Item(6,5).Value() = ' Zdalnie   +1 '

$ws.Cells.Item(7,1).Formula() = '=HYPERLINK("https://nofluffjobs.com/pl/job/remote-automation-qa-with-python-aws-valtech", "https://nofluffjobs.com/pl/job/remote-automation-qa-with-python-aws-valtech")'
$ws.Cells.Item(7,2).Value() = ' Remote Automation QA with Python&AWS  '
$ws.Cells.Item(7,3).Value() = ' Valtech '
$ws.Cells.Item(7,4).Value() = ' 18597  – 27896  PLN'
$ws.Cells.Item(7,5).Value() = ' Zdalnie  '

$ws.Cells.Item(8,1).Formula() = '=HYPERLINK("https://nofluffjobs.com/pl/job/senior-qa-engineer-focal-systems-remote", "https://nofluffjobs.com/pl/job/senior-qa-engineer-focal-systems-remote")'
$ws.Cells.Item(8,2).Value() = ' Senior QA Engineer '
$ws.Cells.Item(8,3).Value() = ' Focal Systems '
$ws.Cells.Item(8,4).Value() = ' 25000  – 35000  PLN'
$ws.Cells.Item(8,5).Value() = ' Zdalnie  '

$ws.Cells.Item(9,1).Formula() = '=HYPERLINK("https://nofluffjobs.com/pl/job/senior-software-tester-macrix-technology-group-remote-1", "https://nofluffjobs.com/pl/job/senior-software-tester-macrix-technology-group-remote-1")'
$ws.Cells.Item(9,2).Value() = ' Senior Software Tester '
$ws.Cells.Item(9,3).Value() = ' Macrix Technology Group '
$ws.Cells.Item(9,4).Value() = ' 14500  – 24000  PLN'
$ws.Cells.Item(9,5).Value() = ' Zdalnie   +1 '

$ws.Cells.Item(10,1).Formula() = '=HYPERLINK("https://nofluffjobs.com/pl/job/mid-senior-network-engineer-with-python-codilime-remote", "https://nofluffjobs.com/pl/job/mid-senior-network-engineer-with-python-codilime-remote")'
$ws.Cells.Item(10,2).Value() = ' Mid/Senior Network Engineer with Python '
$ws.Cells.Item(10,3).Value() = ' CodiLime '
$ws.Cells.Item(10,4).Value() = ' 18000  – 31000  PLN'
$ws.Cells.Item(10,5).Value() = ' Zdalnie   +16 '

$ws.Cells.Item(11,1).Formula() = '=HYPERLINK("https://nofluffjobs.com/pl/job/qa-engineer-link-group-warszawa-2", "https://nofluffjobs.com/pl/job/qa-engineer-link-group-warszawa-2")'
$ws.Cells.Item(11,2).Value() = '  QA Engineer  '
$ws.Cells.Item(11,3).Value() = ' Link Group '
$ws.Cells.Item(11,4).Value() = ' 28000  – 35000  PLN'
$ws.Cells.Item(11,5).Value() = ' Warszawa   +6 '

$ws.Cells.Item(12,1).Formula() = '=HYPERLINK("https://nofluffjobs.com/pl/job/mid-senior-qa-automation-engineer-andea-solutions-remote-1", "https://nofluffjobs.com/pl/job/mid-senior-qa-automation-engineer-andea-solutions-remote-1")'
$ws.Cells.Item(12,2).Value() = ' Mid/Senior QA Automation Engineer '
$ws.Cells.Item(12,3).Value() = ' Andea Solutions Sp. z o.o. '
$ws.Cells.Item(12,4).Value() = ' 12000  – 20000  PLN'
$ws.Cells.Item(12,5).Value() = ' Zdalnie  '

$ws.Cells.Item(13,1).Formula() = '=HYPERLINK("https://nofluffjobs.com/pl/job/qa-engineer-march-networks-poland-gliwice", "https://nofluffjobs.com/pl/job/qa-engineer-march-networks-poland-gliwice")'
$ws.Cells.Item(13,2).Value() = ' QA Engineer '
$ws.Cells.Item(13,3).Value() = ' March Networks Poland Sp. z o.o '
$ws.Cells.Item(13,4).Value() = ' 8000  – 12000  PLN'
$ws.Cells.Item(13,5).Value() = ' Gliwice  '

# Column widths
$ws.Columns.Item(1).ColumnWidth() = 67.28515625
$ws.Columns.Item(2).ColumnWidth() = 46.5703125
$ws.Columns.Item(3).ColumnWidth() = 33.28515625
$ws.Columns.Item(4).ColumnWidth() = 29.7109375
$ws.Columns.Item(5).ColumnWidth() = 22.28515625

# Restore leftover styled-but-empty cell at the new bottom of the extended range
$ws.Cells.Item(210,1).Style() = "Hiperłącze"

# Selection + view
$ws.Range("B25").Select()

# Page setup
$ws.PageSetup.Orientation() = 1
$ws.PageSetup.PaperSize() = 9

